$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old used range (Group/Student table) before writing the new layout
$ws.Cells.Clear()

# New headers
$ws.Range("A1").Value = "Year"
$ws.Range("B1").Value = "EducationName"
$ws.Range("C1").Value = "Average Mark"

# New data rows
$ws.Range("A2").Value = 2019
$ws.Range("B2").Value = "OOP"
$ws.Range("C2").Value = 5

$ws.Range("A3").Value = 2020
$ws.Range("B3").Value = "Math"
$ws.Range("C3").Value = 7

$ws.Range("A4").Value = 2019
$ws.Range("B4").Value = "SUBD"
$ws.Range("C4").Value = 3

$ws.Range("A5").Value = 2020
$ws.Range("B5").Value = "TVIMS"
$ws.Range("C5").Value = 3.5

$ws.Range("A6").Value = 2020
$ws.Range("B6").Value = "OOP"
$ws.Range("C6").Value = 4
